# Columns of the primary types are now cast into those types to ensure
# data integrity on a per-column basis:
#   - Column C  (Special_Rules, text)    : truly-blank cells get "nan"
#   - Column D  (Range, number)          : values stored as booleans are
#                                           re-typed as numbers (same 1/0)
#   - Columns J,K,L,M (Released/Weapon/Armor/Shield, checkbox) : every
#                                           data row is cast to boolean TRUE

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 3
$lastDataRow = 35

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {

    # --- Column C (Special_Rules): blank text cells become the text "nan"
    $cCell = $ws.Cells.Item($r, 3)
    if ($cCell.Value2 -eq "") {
        $cCell.Value = "nan"
    }

    # --- Column D (Range): cells typed as boolean are re-cast to numeric
    $dCell = $ws.Cells.Item($r, 4)
    if ($dCell.Value2.GetType().FullName -eq "System.Boolean") {
        if ($dCell.Value2) {
            $dCell.Value = 1
        } else {
            $dCell.Value = 0
        }
    }

    # --- Columns J, K, L, M (Released/Weapon/Armor/Shield): always cast
    # to a proper boolean TRUE value for every data row.
    $ws.Cells.Item($r, 10).Value = $true
    $ws.Cells.Item($r, 11).Value = $true
    $ws.Cells.Item($r, 12).Value = $true
    $ws.Cells.Item($r, 13).Value = $true
}
